$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Hunk 1: "...the start screen. I created an..." ->
#         "...the start screen class. I created an..."
# Splits the single run into 3 runs:
#   "...start screen" | " class" | ". I created an "
# -----------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("the start screen.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find 'the start screen.'"
}

$insertPos1 = $r1.End - 1  # just before the period
$ip1 = $d.Range($insertPos1, $insertPos1)
$ip1.InsertBefore(" class")
# ip1 now spans the newly-inserted " class" text (Start unchanged, End grew
# to cover it). Toggling Bold on/off on this non-empty range forces a
# persistent run split at both of its edges without leaving any visible
# formatting difference behind (Bold ends up back at its original value).
$ip1.Bold = 1
$ip1.Bold = 0

# -----------------------------------------------------------------
# Hunk 2: "... that made cakes fall to the bottom of the screen. ..." ->
#         "... and an accompany cake class that displayed cakes and
#              made them fall to the bottom of the screen. ..."
# Splits the single run into 3 runs:
#   " " | "and an accompany cake class that displayed cakes and made them" |
#   " fall to the bottom..."
# -----------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("that made cakes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find 'that made cakes'"
}

$oldStart2 = $r2.Start
$oldEnd2 = $r2.End

$ip2 = $d.Range($oldStart2, $oldStart2)
$ip2.InsertBefore("and an accompany cake class that displayed cakes and made them")
# ip2 now spans the inserted replacement text (non-empty range).

# Remove the old "that made cakes" text, which has shifted right by the
# length of the newly-inserted text.
$shift2 = $ip2.End - $ip2.Start
$oldRange2 = $d.Range($oldStart2 + $shift2, $oldEnd2 + $shift2)
$oldRange2.Delete()

# Toggle Bold on the (still valid, non-empty) inserted range to force run
# boundaries on both sides of it, leaving formatting unchanged overall.
$ip2.Bold = 1
$ip2.Bold = 0

Write-Output "Done"
